$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = @(0.0002, 0.0002, 0, 0.1000000014901161, 0)
    3 = @(0.0001, 0.0001, 0, 0.01520000025629997, 0)
    4 = @(0, 0, 0, 0, 0)
    5 = @(0.0141, 0.0141, 0, 0.1773000061511993, 0)
    6 = @(0.0225, 0.0225, 0, 0.2506999969482422, 0)
    7 = @(0.0056, 0.0056, 0, 0.1321000009775162, 0)
    8 = @(0.008, 0.008, 0, 0.1321000009775162, 0)
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Range("B$row").Value = $rowVals[0]
    $ws.Range("C$row").Value = $rowVals[1]
    $ws.Range("D$row").Value = $rowVals[2]
    $ws.Range("E$row").Value = $rowVals[3]
    $ws.Range("F$row").Value = $rowVals[4]
}
